$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4, shifting rows 4..22 down to 5..23
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the new record
$ws.Cells.Item(4, 1).Value = 11
$ws.Cells.Item(4, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(4, 3).Value = "Bíobío"
$ws.Cells.Item(4, 4).Value = 44558
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
$ws.Cells.Item(4, 5).Value = 8
$ws.Cells.Item(4, 6).Value = 100112030
$ws.Cells.Item(4, 7).Value = "Poroto granado"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 250
$ws.Cells.Item(4, 11).Value = 15000
$ws.Cells.Item(4, 12).Value = 16000
$ws.Cells.Item(4, 13).Value = 15400
$ws.Cells.Item(4, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(4, 15).Value = "Región Metropolitana"
$ws.Cells.Item(4, 16).Value = 616
$ws.Cells.Item(4, 17).Value = 25
$ws.Cells.Item(4, 18).Value = "Hortaliza"
